$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header block (rows 2-7): update/re-order the project info rows.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Project Name"
$ws.Range("B2").Value = "OpenCart (Frontend)"

$ws.Range("A3").Value = "Client"
$ws.Range("B3").Value = "OpenCart"

$ws.Range("A4").Value = "Reference Document"
$ws.Range("B4").Value = "FRS"

$ws.Range("A5").Value = "Created By"
$ws.Range("B5").Value = "GM Rajon"

$ws.Range("A6").Value = "Creation Date"
$ws.Range("B6").Value = 45874
$ws.Range("B6").NumberFormat = "mm-dd-yy"

$ws.Range("A7").Value = "Approval Date"
$ws.Range("B7").Value = 45905
$ws.Range("B7").NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------------
# Table header (row 10) - text unchanged, kept for completeness.
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = " Requirement ID"
$ws.Range("B10").Value = "Test Scenario ID"
$ws.Range("C10").Value = "Test Scenario Description"
$ws.Range("D10").Value = "Tes Case ID'S"

# ---------------------------------------------------------------------------
# Row 11 (1.1 Register) - re-centre A11/B11 vertically & force text format on A11.
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = 1.1
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").VerticalAlignment = -4108
$ws.Range("B11").Value = "TS_001"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").VerticalAlignment = -4108
$ws.Range("C11").Value = "Validate the working of Register Account functionality"
$ws.Range("D11").Value = "TC_RF_001`nTC_RF_002`nTC_RF_003`nTC_RF_004`nTC_RF_005`nTC_RF_006`nTC_RF_007`nTC_RF_008`nTC_RF_009`nTC_RF_010`nTC_RF_011`nTC_RF_012`nTC_RF_013`nTC_RF_014`nTC_RF_015`nTC_RF_016`nTC_RF_017`nTC_RF_018`nTC_RF_019`nTC_RF_020`nTC_RF_021`nTC_RF_022`nTC_RF_023`nTC_RF_024`nTC_RF_025`nTC_RF_026`nTC_RF_027`n"
$ws.Rows.Item(11).RowHeight = 352.5

# ---------------------------------------------------------------------------
# Row 12 (1.2 Login)
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = 1.2
$ws.Range("A12").VerticalAlignment = -4108
$ws.Range("B12").Value = "TS_002"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").VerticalAlignment = -4108
$ws.Range("C12").Value = "Validate the working of Login functionality"
$ws.Range("D12").Value = "TC_LF_001`nTC_LF_002`nTC_LF_003`nTC_LF_004`nTC_LF_005`nTC_LF_006`nTC_LF_007`nTC_LF_008`nTC_LF_009`nTC_LF_010`nTC_LF_011`nTC_LF_012`nTC_LF_013`nTC_LF_014`nTC_LF_015`nTC_LF_016`nTC_LF_017`nTC_LF_018`nTC_LF_019`nTC_LF_020`nTC_LF_021`nTC_LF_022`nTC_LF_023`n"
$ws.Rows.Item(12).RowHeight = 302.25

# ---------------------------------------------------------------------------
# Row 13 (1.3 Logout) - new row.
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = 1.3
$ws.Range("A13").HorizontalAlignment = -4108
$ws.Range("A13").VerticalAlignment = -4108
$ws.Range("B13").Value = "TS_003"
$ws.Range("B13").HorizontalAlignment = -4108
$ws.Range("B13").VerticalAlignment = -4108
$ws.Range("C13").Value = "Validate the working of Logout functionality"
$ws.Range("D13").Value = "TC_LF_001`nTC_LF_002`nTC_LF_003`nTC_LF_004`nTC_LF_005`nTC_LF_006`nTC_LF_007`nTC_LF_008`nTC_LF_009`nTC_LF_010`nTC_LF_011`n"
$ws.Range("D13").WrapText = $true
$ws.Range("D13").HorizontalAlignment = -4108
$ws.Range("D13").VerticalAlignment = -4160
$ws.Rows.Item(13).RowHeight = 163.5

# ---------------------------------------------------------------------------
# Row 14 (1.4 Forgot Password) - new row.
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = 1.4
$ws.Range("A14").HorizontalAlignment = -4108
$ws.Range("A14").VerticalAlignment = -4108
$ws.Range("B14").Value = "TS_004"
$ws.Range("B14").HorizontalAlignment = -4108
$ws.Range("B14").VerticalAlignment = -4108
$ws.Range("C14").Value = "Validate the working of Forgot Password functionality"
$ws.Range("C14").WrapText = $true
$ws.Range("D14").Value = "TC_FPF_001`nTC_FPF_002`nTC_FPF_003`nTC_FPF_004`nTC_FPF_005`nTC_FPF_006`nTC_FPF_007`nTC_FPF_008`nTC_FPF_009`nTC_FPF_010`nTC_FPF_011`nTC_FPF_012`nTC_FPF_013`nTC_FPF_014`nTC_FPF_015`nTC_FPF_016`nTC_FPF_017`nTC_FPF_018`nTC_FPF_019`nTC_FPF_020`nTC_FPF_021`nTC_FPF_022`nTC_FPF_023`nTC_FPF_024`nTC_FPF_025`n"
$ws.Range("D14").WrapText = $true
$ws.Range("D14").HorizontalAlignment = -4108
$ws.Range("D14").VerticalAlignment = -4160
$ws.Rows.Item(14).RowHeight = 340.5

# ---------------------------------------------------------------------------
# Row 15 (1.5 Search) - new row.
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = 1.5
$ws.Range("A15").HorizontalAlignment = -4108
$ws.Range("A15").VerticalAlignment = -4108
$ws.Range("B15").Value = "TS_005"
$ws.Range("B15").WrapText = $true
$ws.Range("B15").HorizontalAlignment = -4108
$ws.Range("B15").VerticalAlignment = -4108
$ws.Range("C15").Value = "Validate the working of Seacrch functionality"
$ws.Range("C15").WrapText = $true
$ws.Range("D15").Value = "TC_SF_001`nTC_SF_002`nTC_SF_003`nTC_SF_004`nTC_SF_005`nTC_SF_006`nTC_SF_007`nTC_SF_008`nTC_SF_009`nTC_SF_010`nTC_SF_011`nTC_SF_012`nTC_SF_013`nTC_SF_014`nTC_SF_015`nTC_SF_016`nTC_SF_017`nTC_SF_018`nTC_SF_019`nTC_SF_020`nTC_SF_021`n"
$ws.Range("D15").WrapText = $true
$ws.Range("D15").HorizontalAlignment = -4108
$ws.Range("D15").VerticalAlignment = -4160
$ws.Rows.Item(15).RowHeight = 277.5

# ---------------------------------------------------------------------------
# Row 16 (1.6 Product compare) - new row.
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = 1.6
$ws.Range("A16").HorizontalAlignment = -4108
$ws.Range("A16").VerticalAlignment = -4108
$ws.Range("B16").Value = "TS_006"
$ws.Range("B16").HorizontalAlignment = -4108
$ws.Range("B16").VerticalAlignment = -4108
$ws.Range("C16").Value = "Validate the working of Product compare functionality"
$ws.Range("C16").WrapText = $true
$ws.Range("D16").Value = "TC_PCF_001`nTC_PCF_002`nTC_PCF_003`nTC_PCF_004`nTC_PCF_005`nTC_PCF_006`nTC_PCF_007`nTC_PCF_008`nTC_PCF_009`nTC_PCF_010`nTC_PCF_011`nTC_PCF_012`nTC_PCF_013`nTC_PCF_014`nTC_PCF_015`nTC_PCF_016`nTC_PCF_017`nTC_PCF_018`nTC_PCF_019`nTC_PCF_020`nTC_PCF_021`n"
$ws.Range("D16").WrapText = $true
$ws.Range("D16").HorizontalAlignment = -4108
$ws.Range("D16").VerticalAlignment = -4108
$ws.Rows.Item(16).RowHeight = 277.5

# ---------------------------------------------------------------------------
# Sheet-level view tweaks to mirror the final author session.
# ---------------------------------------------------------------------------
$ws.Range("E16").Select()
$excel.ActiveWindow.ScrollRow = 15
